$d = $word.ActiveDocument
$r = $d.Content

$r.Find.Execute("211×2=422", $true, $false, $false, $false, $false, $true, 1, $false, "545×3=1635", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("891×3=2673", $true, $false, $false, $false, $false, $true, 1, $false, "507×9=4563", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("437×4=1748", $true, $false, $false, $false, $false, $true, 1, $false, "912×8=7296", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("196×9=1764", $true, $false, $false, $false, $false, $true, 1, $false, "332×9=2988", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("377×3=1131", $true, $false, $false, $false, $false, $true, 1, $false, "259×8=2072", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("716×8=5728", $true, $false, $false, $false, $false, $true, 1, $false, "142×2=284", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("626×4=2504", $true, $false, $false, $false, $false, $true, 1, $false, "641×4=2564", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("407×5=2035", $true, $false, $false, $false, $false, $true, 1, $false, "669×6=4014", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("104×5=520", $true, $false, $false, $false, $false, $true, 1, $false, "348×2=696", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("929×8=7432", $true, $false, $false, $false, $false, $true, 1, $false, "158×3=474", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("951×5=4755", $true, $false, $false, $false, $false, $true, 1, $false, "318×8=2544", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("449×7=3143", $true, $false, $false, $false, $false, $true, 1, $false, "540×4=2160", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("779×3=2337", $true, $false, $false, $false, $false, $true, 1, $false, "608×8=4864", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("184×2=368", $true, $false, $false, $false, $false, $true, 1, $false, "713×6=4278", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("497×2=994", $true, $false, $false, $false, $false, $true, 1, $false, "203×9=1827", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("454×4=1816", $true, $false, $false, $false, $false, $true, 1, $false, "624×6=3744", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("501×4=2004", $true, $false, $false, $false, $false, $true, 1, $false, "859×7=6013", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("345×5=1725", $true, $false, $false, $false, $false, $true, 1, $false, "718×3=2154", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("739×8=5912", $true, $false, $false, $false, $false, $true, 1, $false, "709×6=4254", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("308×5=1540", $true, $false, $false, $false, $false, $true, 1, $false, "204×6=1224", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("669×3=2007", $true, $false, $false, $false, $false, $true, 1, $false, "533×4=2132", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("569×7=3983", $true, $false, $false, $false, $false, $true, 1, $false, "152×4=608", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("272×7=1904", $true, $false, $false, $false, $false, $true, 1, $false, "210×3=630", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("466×2=932", $true, $false, $false, $false, $false, $true, 1, $false, "947×5=4735", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("189×3=567", $true, $false, $false, $false, $false, $true, 1, $false, "707×2=1414", 2) | Out-Null
$r = $d.Content
